$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for each data row (2-261).
# All of them currently hold 45178 (2023-09-09) and must be bumped to 45179 (2023-09-10).
$ws.Range("C2:C261").Value = 45179
